$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 89-90: two fixtures swapped id/content ---
# Row 89
$ws.Cells.Item(89, 1).Value = 87
$ws.Cells.Item(89, 2).Value = 6814330
$ws.Cells.Item(89, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(89, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(89, 5).Value = 45137.52083333334
$ws.Cells.Item(89, 6).Value = "NK Maribor"
$ws.Cells.Item(89, 7).Value = "NK Aluminij"
$ws.Cells.Item(89, 8).Value = 1
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = "H"
$ws.Cells.Item(89, 11).Value = 1.363
$ws.Cells.Item(89, 12).Value = 4.5
$ws.Cells.Item(89, 13).Value = 7
$ws.Cells.Item(89, 14).Value = 1.4
$ws.Cells.Item(89, 15).Value = 4.5
$ws.Cells.Item(89, 16).Value = 7
$ws.Cells.Item(89, 17).Value = -1.25
$ws.Cells.Item(89, 18).Value = 1.85
$ws.Cells.Item(89, 19).Value = 1.95
$ws.Cells.Item(89, 20).Value = 2.75
$ws.Cells.Item(89, 21).Value = 1.8
$ws.Cells.Item(89, 22).Value = 2
$ws.Cells.Item(89, 23).Value = 0.3999999999999999
$ws.Cells.Item(89, 24).Value = -1
$ws.Cells.Item(89, 25).Value = -1
$ws.Cells.Item(89, 26).Value = -0.5
$ws.Cells.Item(89, 27).Value = 0.475
$ws.Cells.Item(89, 28).Value = -1
$ws.Cells.Item(89, 29).Value = 1

# Row 90
$ws.Cells.Item(90, 1).Value = 88
$ws.Cells.Item(90, 2).Value = 6814328
$ws.Cells.Item(90, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(90, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(90, 5).Value = 45137.52083333334
$ws.Cells.Item(90, 6).Value = "NK Domzale"
$ws.Cells.Item(90, 7).Value = "NK Bravo"
$ws.Cells.Item(90, 8).Value = 1
$ws.Cells.Item(90, 9).Value = 1
$ws.Cells.Item(90, 10).Value = "D"
$ws.Cells.Item(90, 11).Value = 2.35
$ws.Cells.Item(90, 12).Value = 3.1
$ws.Cells.Item(90, 13).Value = 2.9
$ws.Cells.Item(90, 14).Value = 2.15
$ws.Cells.Item(90, 15).Value = 3.1
$ws.Cells.Item(90, 16).Value = 3.3
$ws.Cells.Item(90, 17).Value = -0.25
$ws.Cells.Item(90, 18).Value = 1.925
$ws.Cells.Item(90, 19).Value = 1.875
$ws.Cells.Item(90, 20).Value = 2.25
$ws.Cells.Item(90, 21).Value = 1.95
$ws.Cells.Item(90, 22).Value = 1.85
$ws.Cells.Item(90, 23).Value = -1
$ws.Cells.Item(90, 24).Value = 2.1
$ws.Cells.Item(90, 25).Value = -1
$ws.Cells.Item(90, 26).Value = -0.5
$ws.Cells.Item(90, 27).Value = 0.4375
$ws.Cells.Item(90, 28).Value = -0.5
$ws.Cells.Item(90, 29).Value = 0.425

# --- Row 185: add result (H/I/J) and AB/AC, update other odds ---
$ws.Cells.Item(185, 1).Value = 183
$ws.Cells.Item(185, 2).Value = 6814408
$ws.Cells.Item(185, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(185, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(185, 5).Value = 45343.45833333334
$ws.Cells.Item(185, 6).Value = "Olimpija Ljubljana"
$ws.Cells.Item(185, 7).Value = "NK Domzale"
$ws.Cells.Item(185, 8).Value = 1
$ws.Cells.Item(185, 9).Value = 0
$ws.Cells.Item(185, 10).Value = "H"
$ws.Cells.Item(185, 11).Value = 1.533
$ws.Cells.Item(185, 12).Value = 4
$ws.Cells.Item(185, 13).Value = 5
$ws.Cells.Item(185, 14).Value = 1.285
$ws.Cells.Item(185, 15).Value = 5.25
$ws.Cells.Item(185, 16).Value = 7.5
$ws.Cells.Item(185, 17).Value = -1.5
$ws.Cells.Item(185, 18).Value = 1.825
$ws.Cells.Item(185, 19).Value = 1.975
$ws.Cells.Item(185, 20).Value = 3
$ws.Cells.Item(185, 21).Value = 1.8
$ws.Cells.Item(185, 22).Value = 2
$ws.Cells.Item(185, 23).Value = 0.2849999999999999
$ws.Cells.Item(185, 24).Value = -1
$ws.Cells.Item(185, 25).Value = -1
$ws.Cells.Item(185, 26).Value = -1
$ws.Cells.Item(185, 27).Value = 0.9750000000000001
$ws.Cells.Item(185, 28).Value = -1
$ws.Cells.Item(185, 29).Value = 1

# --- Rows 186-188: update fixtures (id/date/teams/odds) ---
# Row 186
$ws.Cells.Item(186, 1).Value = 184
$ws.Cells.Item(186, 2).Value = 6816452
$ws.Cells.Item(186, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(186, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(186, 5).Value = 45344.45833333334
$ws.Cells.Item(186, 6).Value = "NK Celje"
$ws.Cells.Item(186, 7).Value = "NK Rogaska"
$ws.Cells.Item(186, 11).Value = 1.25
$ws.Cells.Item(186, 12).Value = 5.75
$ws.Cells.Item(186, 13).Value = 8.5
$ws.Cells.Item(186, 14).Value = 1.285
$ws.Cells.Item(186, 15).Value = 5.25
$ws.Cells.Item(186, 16).Value = 7.5
$ws.Cells.Item(186, 17).Value = -1.5
$ws.Cells.Item(186, 18).Value = 1.825
$ws.Cells.Item(186, 19).Value = 1.975
$ws.Cells.Item(186, 20).Value = 3.25
$ws.Cells.Item(186, 21).Value = 2
$ws.Cells.Item(186, 22).Value = 1.8
$ws.Cells.Item(186, 23).Value = 0
$ws.Cells.Item(186, 24).Value = 0
$ws.Cells.Item(186, 25).Value = 0
$ws.Cells.Item(186, 26).Value = 0
$ws.Cells.Item(186, 27).Value = 0

# Row 187
$ws.Cells.Item(187, 1).Value = 185
$ws.Cells.Item(187, 2).Value = 6814409
$ws.Cells.Item(187, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(187, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(187, 5).Value = 45344.5625
$ws.Cells.Item(187, 6).Value = "NK Maribor"
$ws.Cells.Item(187, 7).Value = "NK Bravo"
$ws.Cells.Item(187, 11).Value = 1.615
$ws.Cells.Item(187, 12).Value = 3.75
$ws.Cells.Item(187, 13).Value = 4.75
$ws.Cells.Item(187, 14).Value = 1.533
$ws.Cells.Item(187, 15).Value = 3.8
$ws.Cells.Item(187, 16).Value = 5.5
$ws.Cells.Item(187, 17).Value = -1
$ws.Cells.Item(187, 18).Value = 1.975
$ws.Cells.Item(187, 19).Value = 1.825
$ws.Cells.Item(187, 20).Value = 2.5
$ws.Cells.Item(187, 21).Value = 1.825
$ws.Cells.Item(187, 22).Value = 1.975
$ws.Cells.Item(187, 23).Value = 0
$ws.Cells.Item(187, 24).Value = 0
$ws.Cells.Item(187, 25).Value = 0
$ws.Cells.Item(187, 26).Value = 0
$ws.Cells.Item(187, 27).Value = 0

# Row 188
$ws.Cells.Item(188, 1).Value = 186
$ws.Cells.Item(188, 2).Value = 6814412
$ws.Cells.Item(188, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(188, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(188, 5).Value = 45346.5625
$ws.Cells.Item(188, 6).Value = "NK Aluminij"
$ws.Cells.Item(188, 7).Value = "Olimpija Ljubljana"
$ws.Cells.Item(188, 11).Value = 5.5
$ws.Cells.Item(188, 12).Value = 4.4
$ws.Cells.Item(188, 13).Value = 1.45
$ws.Cells.Item(188, 14).Value = 5
$ws.Cells.Item(188, 15).Value = 4.333
$ws.Cells.Item(188, 16).Value = 1.5
$ws.Cells.Item(188, 17).Value = 1
$ws.Cells.Item(188, 18).Value = 2
$ws.Cells.Item(188, 19).Value = 1.8
$ws.Cells.Item(188, 20).Value = 2.75
$ws.Cells.Item(188, 21).Value = 1.85
$ws.Cells.Item(188, 22).Value = 1.95
$ws.Cells.Item(188, 23).Value = 0
$ws.Cells.Item(188, 24).Value = 0
$ws.Cells.Item(188, 25).Value = 0
$ws.Cells.Item(188, 26).Value = 0
$ws.Cells.Item(188, 27).Value = 0

# --- Rows 189-192: new fixtures appended ---
# Row 189
$ws.Cells.Item(188, 1).Copy()
$ws.Cells.Item(189, 1).PasteSpecial(-4122)
$ws.Cells.Item(188, 5).Copy()
$ws.Cells.Item(189, 5).PasteSpecial(-4122)
$ws.Cells.Item(189, 1).Value = 187
$ws.Cells.Item(189, 2).Value = 6816451
$ws.Cells.Item(189, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(189, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(189, 5).Value = 45347.375
$ws.Cells.Item(189, 6).Value = "NK Rogaska"
$ws.Cells.Item(189, 7).Value = "NK Maribor"
$ws.Cells.Item(189, 11).Value = 6
$ws.Cells.Item(189, 12).Value = 4.6
$ws.Cells.Item(189, 13).Value = 1.4
$ws.Cells.Item(189, 14).Value = 4.75
$ws.Cells.Item(189, 15).Value = 4.333
$ws.Cells.Item(189, 16).Value = 1.533
$ws.Cells.Item(189, 17).Value = 1
$ws.Cells.Item(189, 18).Value = 1.9
$ws.Cells.Item(189, 19).Value = 1.9
$ws.Cells.Item(189, 20).Value = 3
$ws.Cells.Item(189, 21).Value = 2.025
$ws.Cells.Item(189, 22).Value = 1.775
$ws.Cells.Item(189, 23).Value = 0
$ws.Cells.Item(189, 24).Value = 0
$ws.Cells.Item(189, 25).Value = 0
$ws.Cells.Item(189, 26).Value = 0
$ws.Cells.Item(189, 27).Value = 0

# Row 190
$ws.Cells.Item(188, 1).Copy()
$ws.Cells.Item(190, 1).PasteSpecial(-4122)
$ws.Cells.Item(188, 5).Copy()
$ws.Cells.Item(190, 5).PasteSpecial(-4122)
$ws.Cells.Item(190, 1).Value = 188
$ws.Cells.Item(190, 2).Value = 6814415
$ws.Cells.Item(190, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(190, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(190, 5).Value = 45347.45833333334
$ws.Cells.Item(190, 6).Value = "NK Domzale"
$ws.Cells.Item(190, 7).Value = "NK Celje"
$ws.Cells.Item(190, 11).Value = 6
$ws.Cells.Item(190, 12).Value = 4.6
$ws.Cells.Item(190, 13).Value = 1.4
$ws.Cells.Item(190, 14).Value = 4.75
$ws.Cells.Item(190, 15).Value = 4.2
$ws.Cells.Item(190, 16).Value = 1.55
$ws.Cells.Item(190, 17).Value = 0.75
$ws.Cells.Item(190, 18).Value = 2.1
$ws.Cells.Item(190, 19).Value = 1.7
$ws.Cells.Item(190, 20).Value = 2.75
$ws.Cells.Item(190, 21).Value = 1.8
$ws.Cells.Item(190, 22).Value = 2
$ws.Cells.Item(190, 23).Value = 0
$ws.Cells.Item(190, 24).Value = 0
$ws.Cells.Item(190, 25).Value = 0
$ws.Cells.Item(190, 26).Value = 0
$ws.Cells.Item(190, 27).Value = 0

# Row 191
$ws.Cells.Item(188, 1).Copy()
$ws.Cells.Item(191, 1).PasteSpecial(-4122)
$ws.Cells.Item(188, 5).Copy()
$ws.Cells.Item(191, 5).PasteSpecial(-4122)
$ws.Cells.Item(191, 1).Value = 189
$ws.Cells.Item(191, 2).Value = 6814413
$ws.Cells.Item(191, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(191, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(191, 5).Value = 45347.5625
$ws.Cells.Item(191, 6).Value = "NS Mura"
$ws.Cells.Item(191, 7).Value = "NK Radomlje"
$ws.Cells.Item(191, 11).Value = 1.952
$ws.Cells.Item(191, 12).Value = 3.4
$ws.Cells.Item(191, 13).Value = 3.4
$ws.Cells.Item(191, 14).Value = 2.15
$ws.Cells.Item(191, 15).Value = 3.3
$ws.Cells.Item(191, 16).Value = 3
$ws.Cells.Item(191, 17).Value = -0.25
$ws.Cells.Item(191, 18).Value = 1.95
$ws.Cells.Item(191, 19).Value = 1.85
$ws.Cells.Item(191, 20).Value = 2.5
$ws.Cells.Item(191, 21).Value = 1.975
$ws.Cells.Item(191, 22).Value = 1.825
$ws.Cells.Item(191, 23).Value = 0
$ws.Cells.Item(191, 24).Value = 0
$ws.Cells.Item(191, 25).Value = 0
$ws.Cells.Item(191, 26).Value = 0
$ws.Cells.Item(191, 27).Value = 0

# Row 192
$ws.Cells.Item(188, 1).Copy()
$ws.Cells.Item(192, 1).PasteSpecial(-4122)
$ws.Cells.Item(188, 5).Copy()
$ws.Cells.Item(192, 5).PasteSpecial(-4122)
$ws.Cells.Item(192, 1).Value = 190
$ws.Cells.Item(192, 2).Value = 6814414
$ws.Cells.Item(192, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(192, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(192, 5).Value = 45348.45833333334
$ws.Cells.Item(192, 6).Value = "NK Bravo"
$ws.Cells.Item(192, 7).Value = "FC Koper"
$ws.Cells.Item(192, 11).Value = 2.8
$ws.Cells.Item(192, 12).Value = 3.4
$ws.Cells.Item(192, 13).Value = 2.25
$ws.Cells.Item(192, 14).Value = 2.8
$ws.Cells.Item(192, 15).Value = 3.4
$ws.Cells.Item(192, 16).Value = 2.25
$ws.Cells.Item(192, 17).Value = 0
$ws.Cells.Item(192, 18).Value = 2.1
$ws.Cells.Item(192, 19).Value = 1.7
$ws.Cells.Item(192, 20).Value = 2.5
$ws.Cells.Item(192, 21).Value = 1.975
$ws.Cells.Item(192, 22).Value = 1.825
$ws.Cells.Item(192, 23).Value = 0
$ws.Cells.Item(192, 24).Value = 0
$ws.Cells.Item(192, 25).Value = 0
$ws.Cells.Item(192, 26).Value = 0
$ws.Cells.Item(192, 27).Value = 0

$excel.CutCopyMode = 0